$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("M0 - Account Mgmt")

# Row 4: replace the "Can edit the hash..." note with the new forgotten-password note.
$ws.Range("D4").Value = "Users can't request a forgotten password. Need this as a workaround, or need to configure password reset request."

# Row 6: clear the "Create additional viewmodel first." note (no longer needed).
$ws.Range("D6").ClearContents()

# Row 8: the User ViewModel (with roles) has been added, so this task is done.
$ws.Range("C8").Value = "Done"
$ws.Range("C8").Style = "Good"
$ws.Range("D8").ClearContents()

# Update the active selection to D6, matching the saved view state.
$ws.Activate()
$ws.Range("D6").Select()
